$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.035.38'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '3.603.97'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '242.96'
$ws.Range('E5').Value = '  +2.65%  '
$ws.Range('D6').Value = '657.02'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('D7').Value = '1.68'
$ws.Range('E7').Value = '  +14.69%  '
$ws.Range('D8').Value = '0.417'
$ws.Range('E8').Value = '  +3.55%  '
$ws.Range('E9').Value = '  +6.36%  '
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '3.605.14'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').Value = '43.88'
$ws.Range('E12').Value = '  +3.11%  '
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = '6.44'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '4.276.92'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = '96.805.71'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '0.0000260'
$ws.Range('E17').Value = '  +2.44%  '
$ws.Range('D18').Value = '3.597.60'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '12.80'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '7.76'
$ws.Range('E20').Value = '  -2.21%  '
$ws.Range('D21').Value = '18.03'
$ws.Range('E21').Value = '  +0.76%  '
$ws.Range('D22').Value = '0.534'
$ws.Range('E22').Value = '  +11.49%  '
$ws.Range('D23').Value = '511.42'
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('D24').Value = '3.42'
$ws.Range('E24').Value = '  -3.30%  '
$ws.Range('D25').Value = '0.0000205'
$ws.Range('E25').Value = '  +4.75%  '
$ws.Range('D26').Value = '6.87'
$ws.Range('E26').Value = '  +4.08%  '
$ws.Range('D27').Value = '98.32'
$ws.Range('E27').Value = '  +6.55%  '
$ws.Range('D28').Value = '13.09'
$ws.Range('E28').Value = '  +4.18%  '
$ws.Range('D29').Value = '3.800.79'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').Value = '3.03'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').Value = '0.152'
$ws.Range('E31').Value = '  +9.18%  '
$ws.Range('D32').Value = '11.75'
$ws.Range('E32').Value = '  +4.14%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  +4.64%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').Value = '31.73'
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').Value = '621.00'
$ws.Range('E37').Value = '  +10.79%  '
$ws.Range('D38').Value = '0.571'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('E39').Value = '  +7.78%  '
$ws.Range('D40').Value = '1.62'
$ws.Range('E40').Value = '  +9.48%  '
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('D42').Value = '1.92'
$ws.Range('E42').Value = '  +9.13%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '0.927'
$ws.Range('E44').Value = '  +2.18%  '
$ws.Range('D45').Value = '5.96'
$ws.Range('E45').Value = '  +4.60%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0435'
$ws.Range('E46').Value = '  +5.05%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '2.31'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').Value = '33.44'
$ws.Range('E49').Value = '  -5.40%  '
$ws.Range('D50').Value = '8.49'
$ws.Range('E50').Value = '  +5.21%  '
$ws.Range('E51').Value = '  -0.40%  '
